# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# - "Body" sheet: request body schema collapsed to a single $ref
#   (positionLac.211207Request) instead of inline fields; extra rows removed.
# - "200" sheet: response content schema collapsed to a single $ref
#   (positionLac.211207Response) instead of inline fields; extra rows removed.
# - "400" sheet: response content schema collapsed to a single $ref
#   (errorResponse) instead of inline fields; extra rows removed.
# - "204", "401", "403", "404", "429", "500" sheets: add a new content row
#   referencing the relevant error/response schema ($ref).

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param($ws, $row, $section, $name, $schema)
    $ws.Cells.Item($row, 1).Value = $section
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = "schema"
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = $schema
    $ws.Cells.Item($row, 8).Value = ""
    $ws.Cells.Item($row, 9).Value = "Yes"
    $ws.Cells.Item($row, 10).Value = ""
    $ws.Cells.Item($row, 11).Value = ""
    $ws.Cells.Item($row, 12).Value = ""
    $ws.Cells.Item($row, 13).Value = ""
    $ws.Cells.Item($row, 14).Value = ""
    $ws.Cells.Item($row, 15).Value = ""
}

# --- "Body" sheet: Request body now just references positionLac.211207Request ---
$ws = $wb.Worksheets.Item("Body")
$ws.Rows("4:7").Delete()
Set-SchemaRow $ws 3 "body" "positionLac.211207Request" "positionLac.211207Request"

# --- "200" sheet: Response content now just references positionLac.211207Response ---
$ws = $wb.Worksheets.Item("200")
$ws.Rows("4:9").Delete()
Set-SchemaRow $ws 3 "content" "positionLac.211207Response" "positionLac.211207Response"

# --- "204" sheet: add row referencing positionLac.211207Response ---
$ws = $wb.Worksheets.Item("204")
Set-SchemaRow $ws 3 "content" "positionLac.211207Response" "positionLac.211207Response"

# --- "400" sheet: Response content now just references errorResponse ---
$ws = $wb.Worksheets.Item("400")
$ws.Rows("4:6").Delete()
Set-SchemaRow $ws 3 "content" "errorResponse" "errorResponse"

# --- "401", "403", "404", "429", "500" sheets: add row referencing errorResponse1 ---
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Set-SchemaRow $ws 3 "content" "errorResponse1" "errorResponse1"
}
